$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '63.194.78'
$ws.Range('E2').Value = '  -4.97%  '

$ws.Range('D3').Value = '3.126.70'
$ws.Range('E3').Value = '  -5.58%  '

$ws.Range('E4').Value = '  +0.01%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '564.92'

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '164.14'
$ws.Range('E6').Value = '  -9.29%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.593'
$ws.Range('E7').Value = '  -9.40%  '

$ws.Range('E8').Value = '  +0.02%  '

$ws.Range('D9').Value = '3.125.26'
$ws.Range('E9').Value = '  -5.59%  '

$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.117'
$ws.Range('E10').Value = '  -7.94%  '

$ws.Range('B11').Value = 'Toncoin'
$ws.Range('C11').Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '6.73'
$ws.Range('E11').Value = '  -1.71%  '

$ws.Range('E12').Value = '  -5.77%  '

$ws.Range('D13').Value = '3.678.26'
$ws.Range('E13').Value = '  -5.51%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '0.128'
$ws.Range('E14').Value = '  -1.45%  '

$ws.Range('D15').Value = '63.376.19'
$ws.Range('E15').Value = '  -4.64%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '24.71'
$ws.Range('E16').Value = '  -7.35%  '

$ws.Range('D17').Value = '3.115.91'
$ws.Range('E17').Value = '  -6.27%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.0000155'
$ws.Range('E18').Value = '  -5.79%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '408.53'
$ws.Range('E19').Value = '  -3.61%  '

$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '12.63'
$ws.Range('E20').Value = '  -3.86%  '

$ws.Range('E21').Value = '  -5.27%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '7.03'
$ws.Range('E22').Value = '  -4.13%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '0.999'
$ws.Range('E23').Value = '  -0.12%  '

$ws.Range('E24').Value = '  +0.17%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '68.17'
$ws.Range('E25').Value = '  -4.68%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '0.202'
$ws.Range('E26').Value = '  -1.12%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '0.492'
$ws.Range('E27').Value = '  -4.55%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.0000102'
$ws.Range('E28').Value = '  -11.66%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '8.75'
$ws.Range('E29').Value = '  -4.05%  '

$ws.Range('E30').Value = '  +0.23%  '

$ws.Range('E31').Value = '  -0.07%  '

$ws.Range('E32').Value = '  -6.47%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '21.41'
$ws.Range('E33').Value = '  -4.34%  '

$ws.Range('E34').Value = '  -5.12%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.26'
$ws.Range('E35').Value = '  -5.03%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '155.11'
$ws.Range('E36').Value = '  -3.15%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '1.11'
$ws.Range('E37').Value = '  -6.74%  '

$ws.Range('E38').Value = '  -6.17%  '

$ws.Range('D39').Value = '2.693.42'
$ws.Range('E39').Value = '  -5.77%  '

$ws.Range('E40').Value = '  -7.74%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '4.13'
$ws.Range('E41').Value = '  -4.71%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '23.63'
$ws.Range('E42').Value = '  -10.31%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '38.54'
$ws.Range('E43').Value = '  -2.95%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '0.697'
$ws.Range('E44').Value = '  -7.93%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.0609'
$ws.Range('E45').Value = '  -7.42%  '

$ws.Range('B46').Value = 'VeChain'
$ws.Range('C46').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '0.0257'
$ws.Range('E46').Value = '  -5.81%  '

$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '5.23'
$ws.Range('E47').Value = '  -11.24%  '

$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.99'
$ws.Range('E48').Value = '  -9.22%  '

$ws.Range('B49').Value = 'Bittensor'
$ws.Range('C49').Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '282.01'
$ws.Range('E49').Value = '  -9.23%  '

$ws.Range('E50').Value = '  +0.03%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.0976'
$ws.Range('E51').Value = '  -6.50%  '
